$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values on rows 2, 3, 4, 5, and 15
$ws.Range("F2").Value = -4
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = 2
$ws.Range("F5").Value = -1
$ws.Range("F15").Value = 0
